# Hortaliza, Terminal Hortofrutícola Agro Chillán - Coliflor
# Two new weekly records are inserted at rows 202-203 (pushing the existing
# 202-281 block down to 204-283), matching the new "Fruta / hortaliza,
# semanal" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 202, shifting the rest of the
# table (rows 202:281) down to rows 204:283.
$ws.Rows("202:203").Insert()

# Row 202 - new record (Coliflor, Primera)
$ws.Cells.Item(202, 1).Value = 7
$ws.Cells.Item(202, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(202, 3).Value = "Ñuble"
$ws.Cells.Item(202, 4).Value = 44755
$ws.Cells.Item(202, 5).Value = 16
$ws.Cells.Item(202, 6).Value = 100112008
$ws.Cells.Item(202, 7).Value = "Coliflor"
$ws.Cells.Item(202, 8).Value = "Sin especificar"
$ws.Cells.Item(202, 9).Value = "Primera"
$ws.Cells.Item(202, 10).Value = 200
$ws.Cells.Item(202, 11).Value = 1000
$ws.Cells.Item(202, 12).Value = 1200
$ws.Cells.Item(202, 13).Value = 1100
$ws.Cells.Item(202, 14).Value = "$/unidad"
$ws.Cells.Item(202, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(202, 16).Value = 1100
$ws.Cells.Item(202, 17).Value = 1
$ws.Cells.Item(202, 18).Value = "Hortaliza"

# Row 203 - new record (Coliflor, Segunda)
$ws.Cells.Item(203, 1).Value = 7
$ws.Cells.Item(203, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(203, 3).Value = "Ñuble"
$ws.Cells.Item(203, 4).Value = 44755
$ws.Cells.Item(203, 5).Value = 16
$ws.Cells.Item(203, 6).Value = 100112008
$ws.Cells.Item(203, 7).Value = "Coliflor"
$ws.Cells.Item(203, 8).Value = "Sin especificar"
$ws.Cells.Item(203, 9).Value = "Segunda"
$ws.Cells.Item(203, 10).Value = 150
$ws.Cells.Item(203, 11).Value = 900
$ws.Cells.Item(203, 12).Value = 900
$ws.Cells.Item(203, 13).Value = 900
$ws.Cells.Item(203, 14).Value = "$/unidad"
$ws.Cells.Item(203, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(203, 16).Value = 900
$ws.Cells.Item(203, 17).Value = 1
$ws.Cells.Item(203, 18).Value = "Hortaliza"
